# Apply the textual edits described by the diff.
$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2) | Out-Null
}

# 1. "se le aporta" -> "se aporta"
Replace-Text "se le aporta a la escena" "se aporta a la escena"

# 2. "especificar a que lado" -> "especificar de que lado"
Replace-Text "especificar a que lado" "especificar de que lado"

# 3. "conjunto de cuatro planos" -> "conjunto de seis planos"
Replace-Text "conjunto de cuatro planos" "conjunto de seis planos"

# 4. "A cada habitáculo creado" -> "A cada habitación creada"
Replace-Text "A cada habitáculo creado" "A cada habitación creada"

# 5. "Situamos un plano por cada" -> "Situamos un plano de suelo por cada"
Replace-Text "Situamos un plano por cada" "Situamos un plano de suelo por cada"

# 6. "de los suelos  la podamos" (double space) -> "de los suelos la podamos"
Replace-Text "de los suelos  la podamos" "de los suelos la podamos"
